# Insert a new row at position 182 (shifts existing rows 182:247 down to 183:248)
# and populate the new row with the data for the added Ají (Americana (o)) record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(182).Insert()

$ws.Range("A182").Value = 5
$ws.Range("B182").Value = "Macroferia Regional de Talca"
$ws.Range("C182").Value = "Maule"
$ws.Range("D182").Value = 44704
$ws.Range("E182").Value = 7
$ws.Range("F182").Value = 100112021
$ws.Range("G182").Value = "Ají"
$ws.Range("H182").Value = "Americana (o)"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 50
$ws.Range("K182").Value = 30000
$ws.Range("L182").Value = 30000
$ws.Range("M182").Value = 30000
$ws.Range("N182").Value = "`$/caja 25 kilos"
$ws.Range("O182").Value = "Provincia del Elquí"
$ws.Range("P182").Value = 1200
$ws.Range("Q182").Value = 25
$ws.Range("R182").Value = "Hortaliza"
